$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").Value = '30.694.05'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '1.897.03'
$ws.Range("E3").Value = '  +2.69%  '
$ws.Range("E4").Value = '  +0.05%  '
$__style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'239.12"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = '  +1.34%  '
$__style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.4831"
$ws.Range("D7").Style = $__style
$ws.Range("E7").Value = '  +1.05%  '
$__style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.2850"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = '  +1.72%  '
$__style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.06552"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("D10").Value = '2.001.14'
$ws.Range("E10").Value = '  +8.26%  '
$__style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.07458"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = '  +1.99%  '
$__style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'16.73"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = '  +3.05%  '
$__style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'5.106"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = '  +0.10%  '
$__style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'88.09"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = '  +1.25%  '
$__style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.6672"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = '  +3.54%  '
$ws.Range("D16").Value = '30.679.45'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.282.48'
$ws.Range("E17").Value = '  +8.84%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$__style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'13.31"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$__style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = '  -0.03%  '
$__style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.000007615"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = '  -0.02%  '
$__style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'231.40"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = '  +3.13%  '
$__style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = '  +0.07%  '
$__style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'5.285"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = '  +0.15%  '
$__style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'6.233"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = '  +2.77%  '
$__style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'169.92"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = '  +3.95%  '
$ws.Range("E26").Value = '  +1.54%  '
$__style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'18.77"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = '  +1.69%  '
$__style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'1.968"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = '  +2.90%  '
$__style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'1.405"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = '  -1.57%  '
$__style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.1020"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = '  +11.07%  '
$__style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'4.356"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = '  +2.94%  '
$__style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'4.031"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = '  +2.13%  '
$__style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.05118"
$ws.Range("D33").Style = $__style
$__style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'1.218"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = '  +7.44%  '
$__style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.7580"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("E36").Value = '  +0.59%  '
$__style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.01886"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = '  +4.22%  '
$__style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'2.658"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = '  +2.22%  '
$__style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.9218"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = '  +1.96%  '
$__style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.080"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = '  +1.33%  '
$__style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'107.08"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = '  +0.25%  '
$__style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.4303"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("E43").Value = '  +0.61%  '
$__style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'5.738"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = '  -3.42%  '
$ws.Range("E45").Value = '  +0.80%  '
$__style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'64.66"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = '  +1.19%  '
$__style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.1277"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = '  -2.97%  '
$__style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'1.491"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = '  -4.23%  '
$__style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'8.999"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = '  +3.01%  '
$__style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'33.88"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = '  -0.99%  '
$__style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.05677"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = '  +0.30%  '
